# "Reorganizacion completa: limpieza de modulos antiguos, nuevas entregas y
# optimizacion" - rename the sheet and strip the now-unused header styling
# (bold white-on-blue fill + the 20-char column widths) that the old
# "Datos" template used, leaving plain, unformatted header cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename "Datos" -> "ventas"
$ws.Name = "ventas"

# 2) Remember the header row's text before we touch anything structural.
$headerCount = 6
$headerValues = @()
for ($col = 1; $col -le $headerCount; $col++) {
    $headerValues += $ws.Cells.Item(1, $col).Value()
}

# 3) Drop the explicit 20-character column widths on A:F by deleting and
#    rebuilding those columns, then restore the header text with plain,
#    default formatting (no bold/fill style, no custom width).
$ws.Columns("A:F").Delete()

for ($col = 1; $col -le $headerCount; $col++) {
    $ws.Cells.Item(1, $col).Value = $headerValues[$col - 1]
}

# 4) Make sure the header range carries no leftover style reference
#    (equivalent to the removed bold/centered/blue-fill cell style).
$ws.Range("A1:F1").ClearFormats()
